$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs above the
# existing row 262, so insert a blank row there (shifting 262:368 down to
# 263:369) and fill it in with the new record's values.
$ws.Rows.Item(262).Insert()

$ws.Range("A262").Value = 9
$ws.Range("B262").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C262").Value = 'Metropolitana'
$ws.Range("D262").Value = 45134
$ws.Range("D262").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E262").Value = 13
$ws.Range("F262").Value = 100112026
$ws.Range("G262").Value = 'Haba'
$ws.Range("H262").Value = 'Sin especificar'
$ws.Range("I262").Value = 'Primera'
$ws.Range("J262").Value = 52
$ws.Range("K262").Value = 14000
$ws.Range("L262").Value = 15000
$ws.Range("M262").Value = 14500
$ws.Range("N262").Value = '$/saco 25 kilos'
$ws.Range("O262").Value = 'Provincia de Limarí'
$ws.Range("P262").Value = 580
$ws.Range("Q262").Value = 25
$ws.Range("R262").Value = 'Hortaliza'
